# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (want-to-go count) figures across the four sheets.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 353
$ws.Range("F7").Value = 1159
$ws.Range("F8").Value = 441
$ws.Range("F9").Value = 7086
$ws.Range("F12").Value = 2039
$ws.Range("F13").Value = 7957
$ws.Range("F16").Value = 5500
$ws.Range("F18").Value = 2395
$ws.Range("F19").Value = 1018
$ws.Range("F20").Value = 4559
$ws.Range("F25").Value = 365
$ws.Range("F28").Value = 2297
$ws.Range("F29").Value = 23
$ws.Range("F31").Value = 72
$ws.Range("F32").Value = 134
$ws.Range("F33").Value = 573
$ws.Range("F36").Value = 1478
$ws.Range("F39").Value = 2285
$ws.Range("F40").Value = 2209

# ---- 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 23

# ---- 本地生活 (Local Life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 253
$ws.Range("F3").Value = 1275

# ---- 全部类型 (All Types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 253
$ws.Range("F4").Value = 1275
$ws.Range("F7").Value = 353
$ws.Range("F8").Value = 1159
$ws.Range("F9").Value = 441
$ws.Range("F10").Value = 7087
$ws.Range("F13").Value = 2039
$ws.Range("F14").Value = 7957
$ws.Range("F17").Value = 5500
$ws.Range("F19").Value = 2395
$ws.Range("F20").Value = 1018
$ws.Range("F21").Value = 4559
$ws.Range("F28").Value = 365
$ws.Range("F30").Value = 2297
$ws.Range("F31").Value = 23
$ws.Range("F33").Value = 72
$ws.Range("F34").Value = 134
$ws.Range("F35").Value = 2
$ws.Range("F36").Value = 573
$ws.Range("F39").Value = 23
$ws.Range("F40").Value = 1478
$ws.Range("F43").Value = 2286
$ws.Range("F45").Value = 2209
